$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 279, shifting existing rows 279:366 down to 280:367
$ws.Rows("279:279").Insert()

# Populate the newly inserted row 279 with the new data point.
# Columns A, B, C, E, F, G, H, I, Q, R are identical across the dataset,
# so copy them from the (now shifted) following row 280.
$ws.Range("A279").Value = $ws.Range("A280").Value()
$ws.Range("B279").Value = $ws.Range("B280").Value()
$ws.Range("C279").Value = $ws.Range("C280").Value()
$ws.Range("D279").Value = 44985
$ws.Range("E279").Value = $ws.Range("E280").Value()
$ws.Range("F279").Value = $ws.Range("F280").Value()
$ws.Range("G279").Value = $ws.Range("G280").Value()
$ws.Range("H279").Value = $ws.Range("H280").Value()
$ws.Range("I279").Value = $ws.Range("I280").Value()
$ws.Range("J279").Value = 70
$ws.Range("K279").Value = 30000
$ws.Range("L279").Value = 32000
$ws.Range("M279").Value = 31000
$ws.Range("N279").Value = "$/saco 25 kilos"
$ws.Range("O279").Value = "Región Metropolitana"
$ws.Range("P279").Value = 1240
$ws.Range("Q279").Value = $ws.Range("Q280").Value()
$ws.Range("R279").Value = $ws.Range("R280").Value()
